# Populate the four previously-empty "main page" tables (sheets 1, 3, 4, 5)
# with the rows/columns/year/value tables described by the commit
# ("integrate code for date top issue and main page more than 1 table").
# Sheets 2 stays empty; sheets 6-8 already carried this data and are
# left untouched.

$wb = $excel.ActiveWorkbook

function Fill-Table($ws, $rows) {
    # Header row, always the same 4 labels, bold + bordered + centered
    # (matches the pre-existing header style used on the other sheets).
    $ws.Cells.Item(1,1).Value = "rows"
    $ws.Cells.Item(1,2).Value = "columns"
    $ws.Cells.Item(1,3).Value = "year"
    $ws.Cells.Item(1,4).Value = "value"

    $headerRange = $ws.Range("A1:D1")
    $headerRange.Font.Bold = $true
    $headerRange.HorizontalAlignment = -4108
    $headerRange.VerticalAlignment = -4160
    $headerRange.Borders.LineStyle = 1

    $r = 2
    foreach ($row in $rows) {
        $ws.Cells.Item($r,1).Value = $row[0]
        if ($null -ne $row[1]) {
            $ws.Cells.Item($r,2).Value = $row[1]
        }
        $ws.Cells.Item($r,3).Value = $row[2]
        $ws.Cells.Item($r,4).Value = $row[3]
        $r = $r + 1
    }
}

# --- Sheet 1: "10__bc529b88-4ffa-36" -> 10 Inventories ---------------------
$ws1 = $wb.Worksheets.Item(1)
$sheet1Rows = @(
    ,@("10 Inventories Finished goods", $null, 2022, 21156)
    ,@("10 Inventories Goods in transit", $null, 2022, 13236)
    ,@("10 Inventories Provision for inventory obsolescence", $null, 2022, -73)
    ,@("10 Inventories ", $null, 2022, 34319)
    ,@("10 Inventories Finished goods", $null, 2021, 14472)
    ,@("10 Inventories Goods in transit", $null, 2021, 13324)
    ,@("10 Inventories Provision for inventory obsolescence", $null, 2021, -161)
    ,@("10 Inventories ", $null, 2021, 27635)
)
Fill-Table $ws1 $sheet1Rows

# --- Sheet 3: "12__7f5b7384-15a4-38" -> Trade and other payables -----------
$ws3 = $wb.Worksheets.Item(3)
$sheet3Rows = @(
    ,@(" Trade payables", "Trade and other payables", 2022, 1564)
    ,@(" Amounts payable to controlling entity (Refer Note 19)", "Trade and other payables", 2022, 19828)
    ,@(" Other creditors and accruals", "Trade and other payables", 2022, 15277)
    ,@(" ", "Trade and other payables", 2022, 36669)
    ,@(" Trade payables", "$'000", 2021, 3752)
    ,@(" Amounts payable to controlling entity (Refer Note 19)", "$'000", 2021, 19220)
    ,@(" Other creditors and accruals", "$'000", 2021, 12092)
    ,@(" ", "$'000", 2021, 35064)
)
Fill-Table $ws3 $sheet3Rows

# --- Sheet 4: "13__7f5b7384-15a4-38" -> 13 Lease liabilities ---------------
$ws4 = $wb.Worksheets.Item(4)
$sheet4Rows = @(
    ,@("13 Lease liabilities Opening balance", $null, 2022, 11043)
    ,@("13 Lease liabilities Additions", $null, 2022, 888)
    ,@("13 Lease liabilities Payments", $null, 2022, -4452)
    ,@("13 Lease liabilities Closing balance", $null, 2022, 7638)
    ,@("13 Lease liabilities Current", $null, 2022, 3652)
    ,@("13 Lease liabilities Non-current", $null, 2022, 3986)
    ,@("Amounts recognised in profit or loss Interest on lease liabilities", $null, 2022, 159)
    ,@("Amounts recognised in profit or loss Depreciation of right-of use assets", $null, 2022, 2649)
    ,@("Amounts recognised in the statement of cash flows Total cash outflow for leases", $null, 2022, 4293)
    ,@("13 Lease liabilities Opening balance", $null, 2021, 1289)
    ,@("13 Lease liabilities Additions", $null, 2021, 11685)
    ,@("13 Lease liabilities Payments", $null, 2021, -2009)
    ,@("13 Lease liabilities Closing balance", $null, 2021, 11.043)
    ,@("13 Lease liabilities Current", $null, 2021, 4090)
    ,@("13 Lease liabilities Non-current", $null, 2021, 6953)
    ,@("Amounts recognised in profit or loss Interest on lease liabilities", $null, 2021, 78)
    ,@("Amounts recognised in profit or loss Depreciation of right-of use assets", $null, 2021, 1.517)
    ,@("Amounts recognised in the statement of cash flows Total cash outflow for leases", $null, 2021, 1931)
)
Fill-Table $ws4 $sheet4Rows

# --- Sheet 5: "14__7f5b7384-15a4-38" -> Current / Non-current employee leave
$ws5 = $wb.Worksheets.Item(5)
$sheet5Rows = @(
    ,@("Current Liability for annual leave and other current employee benefits", $null, 2022, 1.139)
    ,@("Current ", $null, 2022, 1139)
    ,@("Non-current Liability for long-service leave", $null, 2022, 787)
    ,@("Non-current ", $null, 2022, 787)
    ,@("Current Liability for annual leave and other current employee benefits", $null, 2021, 1218)
    ,@("Current ", $null, 2021, 1218)
    ,@("Non-current Liability for long-service leave", $null, 2021, 720)
    ,@("Non-current ", $null, 2021, 720)
)
Fill-Table $ws5 $sheet5Rows

Write-Host "Populated sheets 1, 3, 4 and 5 with the extracted note tables."
